$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.283.45'
$ws.Range('E2').Value = '  +3.88%  '
$ws.Range('D3').Value = '2.428.23'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.98'
$ws.Range('E5').Value = '  +3.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.41'
$ws.Range('E6').Value = '  +5.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.517'
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').Value = '  +7.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.65'
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.16'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.06'
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').Value = '2.807.32'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '2.457.52'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.846'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').Value = '45.173.56'
$ws.Range('E18').Value = '  +3.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.25'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.37'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').Value = '0.0₃0920'
$ws.Range('E21').Value = '  +2.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.81'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '244.20'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.59'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('B30').Value = 'OKB'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '49.09'
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.94'
$ws.Range('E31').Value = '  +1.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.37'
$ws.Range('E32').Value = '  +10.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.126'
$ws.Range('E33').Value = '  +4.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.22'
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.47'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.86'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '125.47'
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.73'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0290'
$ws.Range('E44').Value = '  +2.13%  '
$ws.Range('D45').Value = '1.939.80'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').Value = '  +3.49%  '
$ws.Range('E48').Value = '  +16.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.14'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '76.50'
$ws.Range('E50').Value = '  +5.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.03'
$ws.Range('E51').Value = '  +2.61%  '
